function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").Value = "28.273.42"
    $ws.Range("E2").Value = "  -0.45%  "
    $ws.Range("D3").Value = "1.805.31"
    $ws.Range("E3").Value = "  -0.99%  "
    Set-TextValue $ws "D4" "1.006"
    $ws.Range("E4").Value = "  +0.42%  "
    Set-TextValue $ws "D5" "314.78"
    $ws.Range("E5").Value = "  -0.08%  "
    Set-TextValue $ws "D6" "1.005"
    $ws.Range("E6").Value = "  +0.34%  "
    Set-TextValue $ws "D7" "0.5231"
    $ws.Range("E7").Value = "  +1.86%  "
    Set-TextValue $ws "D8" "0.3825"
    $ws.Range("E8").Value = "  -2.58%  "
    Set-TextValue $ws "D9" "0.07944"
    $ws.Range("E9").Value = "  +3.61%  "
    Set-TextValue $ws "D10" "41.88"
    $ws.Range("E10").Value = "  +0.66%  "
    Set-TextValue $ws "D11" "1.099"
    $ws.Range("E11").Value = "  -0.88%  "
    Set-TextValue $ws "D12" "6.327"
    $ws.Range("E12").Value = "  +0.91%  "
    $ws.Range("E13").Value = "  +0.46%  "
    Set-TextValue $ws "D14" "20.67"
    $ws.Range("E14").Value = "  -1.61%  "
    $ws.Range("D15").Value = "1.810.33"
    $ws.Range("E15").Value = "  -0.82%  "
    Set-TextValue $ws "D16" "7.345"
    $ws.Range("E16").Value = "  -1.93%  "
    Set-TextValue $ws "D17" "92.71"
    $ws.Range("E17").Value = "  -0.65%  "
    Set-TextValue $ws "D18" "0.00001093"
    $ws.Range("E18").Value = "  -0.21%  "
    Set-TextValue $ws "D19" "0.06600"
    $ws.Range("E19").Value = "  -0.96%  "
    Set-TextValue $ws "D20" "1.005"
    $ws.Range("E20").Value = "  +0.39%  "
    Set-TextValue $ws "D21" "17.42"
    $ws.Range("E21").Value = "  -1.54%  "
    Set-TextValue $ws "D22" "5.970"
    $ws.Range("E22").Value = "  -2.49%  "
    $ws.Range("D23").Value = "28.316.05"
    $ws.Range("E23").Value = "  -0.36%  "
    Set-TextValue $ws "D24" "11.18"
    $ws.Range("E24").Value = "  +0.03%  "
    Set-TextValue $ws "D25" "2.249"
    $ws.Range("E25").Value = "  -0.20%  "
    Set-TextValue $ws "D26" "157.57"
    $ws.Range("E26").Value = "  +0.73%  "
    $ws.Range("E27").Value = "  -1.23%  "
    $ws.Range("D28").Value = "2.016.06"
    $ws.Range("E28").Value = "  -0.91%  "
    Set-TextValue $ws "D29" "2.395"
    $ws.Range("E29").Value = "  +0.13%  "
    Set-TextValue $ws "D30" "123.23"
    $ws.Range("E30").Value = "  -0.75%  "
    Set-TextValue $ws "D31" "0.1105"
    $ws.Range("E31").Value = "  +1.34%  "
    $ws.Range("E32").Value = "  -4.36%  "
    Set-TextValue $ws "D33" "3.669"
    $ws.Range("E33").Value = "  +0.32%  "
    Set-TextValue $ws "D34" "5.587"
    $ws.Range("E34").Value = "  -1.11%  "
    Set-TextValue $ws "D35" "0.07205"
    $ws.Range("E35").Value = "  +1.19%  "
    Set-TextValue $ws "D36" "12.14"
    $ws.Range("E36").Value = "  +8.29%  "
    Set-TextValue $ws "D37" "0.2171"
    $ws.Range("E37").Value = "  -1.78%  "
    Set-TextValue $ws "D38" "0.02315"
    $ws.Range("E38").Value = "  -0.49%  "
    Set-TextValue $ws "D39" "8.713"
    $ws.Range("E39").Value = "  -0.86%  "
    Set-TextValue $ws "D40" "5.044"
    $ws.Range("E40").Value = "  -2.40%  "
    Set-TextValue $ws "D41" "0.6192"
    $ws.Range("E41").Value = "  -1.03%  "
    Set-TextValue $ws "D42" "1.168"
    $ws.Range("E42").Value = "  +0.06%  "
    Set-TextValue $ws "D43" "1.381"
    $ws.Range("E43").Value = "  -0.76%  "
    Set-TextValue $ws "D44" "0.6045"
    $ws.Range("E44").Value = "  +2.77%  "
    $ws.Range("B45").Value = "EnergySwap"
    $ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    Set-TextValue $ws "D45" "13.23"
    $ws.Range("E45").Value = "  -1.44%  "
    $ws.Range("B46").Value = "PancakeSwap"
    $ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
    Set-TextValue $ws "D46" "3.775"
    $ws.Range("E46").Value = "  +1.60%  "
    Set-TextValue $ws "D47" "125.60"
    $ws.Range("E47").Value = "  +0.78%  "
    Set-TextValue $ws "D48" "1.206"
    $ws.Range("E48").Value = "  +0.94%  "
    $ws.Range("E49").Value = "  -2.37%  "
    Set-TextValue $ws "D50" "0.06829"
    $ws.Range("E50").Value = "  -0.99%  "
    Set-TextValue $ws "D51" "73.09"
    $ws.Range("E51").Value = "  -1.12%  "
